$d = $word.ActiveDocument

# The last paragraph in the document is the final green bullet line
# ("- Print the validity result and message returned from
# validate_password for each tested password"). Insert six new
# green-colored paragraphs after it, matching its font color exactly.

$anchor = $d.Paragraphs.Last.Range
$color = $anchor.Font.Color

$newLines = @(
    '1. Password must be at least 8 characters long.',
    '2. Password must include at least one number.',
    '3. Password must include at least one special character from the set: !@#$%^&*(),.?":{}|<>.',
    '4. Address must be at least 8 characters long.',
    '5. Address must include at least one number.',
    '6. Address must include at least one special character from the set: !@#$%^&*(),.?":{}|<>.'
)

$current = $anchor
foreach ($line in $newLines) {
    $current.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last.Range
    $newPara.InsertAfter($line)
    $textOnly = $d.Range($newPara.Start, $newPara.End - 1)
    $textOnly.Font.Color = $color
    $current = $newPara
}
